# Adds the four new "UsersN" worksheets (Users2..Users5) that mirror the
# existing Users / Users1 sheets, reflecting a bugfixed scraper run that
# picked up more links per user and an additional scraped user (PR).
$wb = $excel.ActiveWorkbook

function Add-UsersSheet {
    param(
        [string]$SheetName,
        [object[]]$Rows
    )

    # New sheet goes at the end of the tab strip, after the current last sheet.
    $last = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $last)
    $ws.Name = $SheetName

    # Header row (columns B:F) - column A has no header, matching the source sheets.
    $ws.Cells.Item(1, 2).Value = "name"
    $ws.Cells.Item(1, 3).Value = "email"
    $ws.Cells.Item(1, 4).Value = "links"
    $ws.Cells.Item(1, 5).Value = "buzzwords"
    $ws.Cells.Item(1, 6).Value = "superbuzzwords"

    # Data rows (columns A:F), starting at row 2.
    $r = 2
    foreach ($row in $Rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $ws.Cells.Item($r, 5).Value = $row[4]
        $ws.Cells.Item($r, 6).Value = $row[5]
        $r++
    }
}


$Users2 = @(
    @(0, 'KiP', 'kirsten.preis@amst.at', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://www.reiser-st.com/'': False, ''https://www.amst.co.at/en/aerospace-medicine/'': False, ''https://www.amst.co.at/en/civil-aviation/'': False}', '[''Airbus'', ''Simulator'', ''Lufthansa'', ''Training'', ''Flugschule'', ''VR'', ''XR'', ''EASA'', ''Pilatus'', ''Simulation'', ''Brunner'', ''Upset'', ''UPRT'']', '[''Simulator'', ''Training'', ''Flugschule'']'),
    @(1, 'Kiki', 'whizzogalaxy@web.de', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.flighttraining-service.de/'': False}', '[''Flugschule'', ''Jobs'', ''Training'', ''PPL'', ''CPL'', ''Bayern'', ''Salzburg'', ''Österreich'']', '[]'),
    @(2, 'Kiki', 'kirsten.preis@flightteam.de', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://air-munich.de/'': False, ''http://www.fliegerverein.eu/'': False, ''https://www.mfa.aero/de/'': False, ''https://www.flugausbildung.de/'': False, ''https://www.eaa.aero/en/'': False}', '[''Training, Simulator, PPL, UL, Lehrgang, ATPL, CPL'']', '[''Flightteam'', ''reise'']'),
    @(3, 'PR', 'peter@rothweb.at', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://air-munich.de/'': False, ''http://www.fliegerverein.eu/'': False, ''https://www.mfa.aero/de/'': False, ''https://www.flugausbildung.de/'': False, ''https://www.eaa.aero/en/'': False, ''https://www.reiser-st.com/'': False, ''https://www.amst.co.at/en/aerospace-medicine/'': False, ''https://www.amst.co.at/en/civil-aviation/'': False}', '[''VR'', ''XR'', ''unity'', ''varjo'', ''simulation'', ''simulator'']', '[]')
)
Add-UsersSheet "Users2" $Users2


$Users3 = @(
    @(0, 'KiP', 'kirsten.preis@amst.at', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://www.reiser-st.com/'': False, ''https://www.amst.co.at/en/aerospace-medicine/'': False, ''https://www.amst.co.at/en/civil-aviation/'': False}', '[''Airbus'', ''Simulator'', ''Lufthansa'', ''Training'', ''Flugschule'', ''VR'', ''XR'', ''EASA'', ''Pilatus'', ''Simulation'', ''Brunner'', ''Upset'', ''UPRT'']', '[''Simulator'', ''Training'', ''Flugschule'']'),
    @(1, 'Kiki', 'whizzogalaxy@web.de', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.flighttraining-service.de/'': False}', '[''Flugschule'', ''Jobs'', ''Training'', ''PPL'', ''CPL'', ''Bayern'', ''Salzburg'', ''Österreich'']', '[]'),
    @(2, 'Kiki', 'kirsten.preis@flightteam.de', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://air-munich.de/'': False, ''http://www.fliegerverein.eu/'': False, ''https://www.mfa.aero/de/'': False, ''https://www.flugausbildung.de/'': False, ''https://www.eaa.aero/en/'': False}', '[''Training, Simulator, PPL, UL, Lehrgang, ATPL, CPL'']', '[''Flightteam'', ''reise'']'),
    @(3, 'PR', 'peter@rothweb.at', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://air-munich.de/'': False, ''http://www.fliegerverein.eu/'': False, ''https://www.mfa.aero/de/'': False, ''https://www.flugausbildung.de/'': False, ''https://www.eaa.aero/en/'': False, ''https://www.reiser-st.com/'': False, ''https://www.amst.co.at/en/aerospace-medicine/'': False, ''https://www.amst.co.at/en/civil-aviation/'': False}', '[''VR'', ''XR'', ''unity'', ''varjo'', ''simulation'', ''simulator'']', '[]'),
    @(4, 'PR', 'peter@rothweb.at', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://www.reiser-st.com/'': False, ''https://www.amst.co.at/en/aerospace-medicine/'': False, ''https://www.amst.co.at/en/civil-aviation/'': False}', '[''VR'', ''XR'', ''unity'', ''varjo'', ''simulation'', ''simulator'']', '[]')
)
Add-UsersSheet "Users3" $Users3


$Users4 = @(
    @(0, 'KiP', 'kirsten.preis@amst.at', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://www.reiser-st.com/'': False, ''https://www.amst.co.at/en/aerospace-medicine/'': False, ''https://www.amst.co.at/en/civil-aviation/'': False}', '[''Airbus'', ''Simulator'', ''Lufthansa'', ''Training'', ''Flugschule'', ''VR'', ''XR'', ''EASA'', ''Pilatus'', ''Simulation'', ''Brunner'', ''Upset'', ''UPRT'']', '[''Simulator'', ''Training'', ''Flugschule'']'),
    @(1, 'Kiki', 'whizzogalaxy@web.de', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.flighttraining-service.de/'': False}', '[''Flugschule'', ''Jobs'', ''Training'', ''PPL'', ''CPL'', ''Bayern'', ''Salzburg'', ''Österreich'']', '[]'),
    @(2, 'Kiki', 'kirsten.preis@flightteam.de', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://air-munich.de/'': False, ''http://www.fliegerverein.eu/'': False, ''https://www.mfa.aero/de/'': False, ''https://www.flugausbildung.de/'': False, ''https://www.eaa.aero/en/'': False}', '[''Training, Simulator, PPL, UL, Lehrgang, ATPL, CPL'']', '[''Flightteam'', ''reise'']'),
    @(3, 'PR', 'peter@rothweb.at', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://www.reiser-st.com/'': False, ''https://www.amst.co.at/en/aerospace-medicine/'': False, ''https://www.amst.co.at/en/civil-aviation/'': False}', '[''VR'', ''XR'', ''unity'', ''varjo'', ''simulation'', ''simulator'']', '[]')
)
Add-UsersSheet "Users4" $Users4


$Users5 = @(
    @(0, 'KiP', 'kirsten.preis@amst.at', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://www.reiser-st.com/'': False, ''https://www.amst.co.at/en/aerospace-medicine/'': False, ''https://www.amst.co.at/en/civil-aviation/'': False}', '[''Airbus'', ''Simulator'', ''Lufthansa'', ''Training'', ''Flugschule'', ''VR'', ''XR'', ''EASA'', ''Pilatus'', ''Simulation'', ''Brunner'', ''Upset'', ''UPRT'']', '[''Simulator'', ''Training'', ''Flugschule'']'),
    @(1, 'Kiki', 'whizzogalaxy@web.de', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.flighttraining-service.de/'': False}', '[''Flugschule'', ''Jobs'', ''Training'', ''PPL'', ''CPL'', ''Bayern'', ''Salzburg'', ''Österreich'']', '[]'),
    @(2, 'Kiki', 'kirsten.preis@flightteam.de', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://air-munich.de/'': False, ''http://www.fliegerverein.eu/'': False, ''https://www.mfa.aero/de/'': False, ''https://www.flugausbildung.de/'': False, ''https://www.eaa.aero/en/'': False}', '[''Training, Simulator, PPL, UL, Lehrgang, ATPL, CPL'']', '[''Flightteam'', ''reise'']'),
    @(3, 'PR', 'peter@rothweb.at', '{''https://www.flugrevue.de/'': True, ''https://www.aero.de/'': True, ''https://www.pressebox.de/'': True, ''https://www.etcusa.com/'': False, ''https://www.flighttraining-service.de/'': False, ''https://www.reiser-st.com/'': False, ''https://www.amst.co.at/en/aerospace-medicine/'': False, ''https://www.amst.co.at/en/civil-aviation/'': False}', '[''VR'', ''XR'', ''unity'', ''varjo'', ''simulation'', ''simulator'']', '[]')
)
Add-UsersSheet "Users5" $Users5
